$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.909.32"
$ws.Range("E2").Value = "  +1.50%  "
# Row 3
$ws.Range("D3").Value = "1.776.58"
$ws.Range("E3").Value = "  +1.64%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4548"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.92%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3593"
$ws.Range("D8").Style = "Normal"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.22%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.01%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.39%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.16%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.74%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.053"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.70%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.223"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.73%  "
# Row 16
$ws.Range("D16").Value = "1.774.24"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.95%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001065"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06443"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.67%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.72%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.808"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "
# Row 23
$ws.Range("D23").Value = "27.958.73"
$ws.Range("E23").Value = "  +1.46%  "
# Row 24
$ws.Range("E24").Value = "  +1.35%  "
# Row 25
$ws.Range("E25").Value = "  -0.22%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.16%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.75%  "
# Row 28
$ws.Range("D28").Value = "1.979.91"
$ws.Range("E28").Value = "  +1.49%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.213"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.50%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.20%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.122"
$ws.Range("D31").Style = "Normal"
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09226"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.00%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.667"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.573"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.04%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.12%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02302"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.11%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06193"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.44%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2104"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "
# Row 39
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.996"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.25%  "
# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6341"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.13%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.188"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.80%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.388"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.41%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.913"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.24%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.51%  "
# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5925"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "
# Row 46
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.739"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.56%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.19%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.962"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.87%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06931"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.34%  "
# Row 50
$ws.Range("E50").Value = "  -0.16%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.47%  "
